# Updates "phan chia cong viec" (task-division) sheet:
#  - Khanh now covers both PM + BA-analysis roles (rows 4-5 merged)
#  - Khoa moves from the removed "Deployment" role onto Designer (new row 8)
#  - Every other member shifts up one row, STT renumbered 1-10
#  - trailing spacer row 15 added below the table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4/5: Khanh, PM + BA-analysis, merged vertically ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Khanh"
$ws.Range("C4").Value = "Trưởng nhóm (PM)"
$ws.Range("C5").Value = "BA – Phân tích nghiệp vụ"

# clear the pre-existing box borders before merging/re-drawing them,
# so the merge does not inherit a stray fill color on the split edges
$ws.Range("A4:B5").Borders.LineStyle = -4142

[void]$ws.Range("A4:A5").Merge()
[void]$ws.Range("B4:B5").Merge()

$ws.Range("A4").Borders.Item(7).LineStyle = 1
$ws.Range("A4").Borders.Item(10).LineStyle = 1
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A5").Borders.Item(7).LineStyle = 1
$ws.Range("A5").Borders.Item(10).LineStyle = 1
$ws.Range("A5").Borders.Item(9).LineStyle = 1

$ws.Range("B4").Borders.Item(7).LineStyle = 1
$ws.Range("B4").Borders.Item(10).LineStyle = 1
$ws.Range("B4").Borders.Item(8).LineStyle = 1
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B5").Borders.Item(7).LineStyle = 1
$ws.Range("B5").Borders.Item(10).LineStyle = 1
$ws.Range("B5").Borders.Item(9).LineStyle = 1
$ws.Range("B5").HorizontalAlignment = -4108

# --- Rows 6-14: renumbered / reshuffled member table ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Ngôn"
$ws.Range("C6").Value = "BA – Viết tài liệu yêu cầu"

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Kiệt"
$ws.Range("C7").Value = "Designer – UI/UX"

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Khoa"
$ws.Range("C8").Value = "Designer – UI/UX"

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Gia Huy"
$ws.Range("C9").Value = "Developer Backend"

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "Hòa Hộp"
$ws.Range("C10").Value = "Developer Backend"

$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Nghĩa"
$ws.Range("C11").Value = "Developer Frontend"

$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Vinh Huy"
$ws.Range("C12").Value = "Developer Frontend"

$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Khải"
$ws.Range("C13").Value = "Tester – Kiểm thử"

$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Phát"
$ws.Range("C14").Value = "Tester – Kiểm thử"

# row 13/14: give the bottom "Tester" pair the same broken-box look as
# the PM/BA block above (C13 loses its bottom edge, B14 loses its right edge)
$ws.Range("C13").Borders.LineStyle = -4142
$ws.Range("C13").Borders.Item(7).LineStyle = 1
$ws.Range("C13").Borders.Item(10).LineStyle = 1
$ws.Range("C13").Borders.Item(8).LineStyle = 1

$ws.Range("B14").Borders.Item(7).LineStyle = 1
$ws.Range("B14").Borders.Item(8).LineStyle = 1
$ws.Range("B14").Borders.Item(9).LineStyle = 1

# --- Row 15: empty trailing spacer cell, styled like the role column ---
$ws.Rows.Item(15).RowHeight = 24.9
$c15 = $ws.Range("C15")
$c15.Font.Name = "Times New Roman"
$c15.Font.Size = 11
$c15.Font.Color = 0
$c15.WrapText = $true
$c15.VerticalAlignment = -4108

[void]$ws.Range("D5").Select()
